$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws4 = $wb.Worksheets.Item(4)

# Sheet1 (展览) simple F-value updates
$ws1.Range("F2").Value2 = 20852
$ws1.Range("F3").Value2 = 810
$ws1.Range("F4").Value2 = 335
$ws1.Range("F5").Value2 = 1118
$ws1.Range("F6").Value2 = 22
$ws1.Range("F7").Value2 = 7796
$ws1.Range("F8").Value2 = 544
$ws1.Range("F9").Value2 = 20
$ws1.Range("F10").Value2 = 753
$ws1.Range("F11").Value2 = 296
$ws1.Range("F12").Value2 = 52
$ws1.Range("F13").Value2 = 180
$ws1.Range("F14").Value2 = 149
$ws1.Range("F15").Value2 = 25
$ws1.Range("F18").Value2 = 1353
$ws1.Range("F19").Value2 = 500
$ws1.Range("F20").Value2 = 78
$ws1.Range("F23").Value2 = 77
$ws1.Range("F24").Value2 = 81
$ws1.Range("F25").Value2 = 341
$ws1.Range("F26").Value2 = 1161
$ws1.Range("F28").Value2 = 33
$ws1.Range("F29").Value2 = 207
$ws1.Range("F30").Value2 = 5213
$ws1.Range("F31").Value2 = 590
$ws1.Range("F32").Value2 = 114
$ws1.Range("F37").Value2 = 58
$ws1.Range("F38").Value2 = 12933
$ws1.Range("F39").Value2 = 1352
$ws1.Range("F40").Value2 = 110
$ws1.Range("F41").Value2 = 43
$ws1.Range("F42").Value2 = 65
$ws1.Range("F43").Value2 = 293
$ws1.Range("F44").Value2 = 407
$ws1.Range("F45").Value2 = 4036

# Sheet1 rows 33-36 rotation updates
$ws1.Range("C33").Value2 = "苏州·第三届.OCG.Summer Carnival-国潮动漫游戏嘉年华"
$ws1.Range("E33").Value2 = "2024.08.10 09:00-08.11 17:00"
$ws1.Range("F33").Value2 = 4970
$ws1.Range("G33").Value2 = 70
$ws1.Range("H33").Value2 = "https://show.bilibili.com/platform/detail.html?id=88451"
$ws1.Range("I33").Value2 = "//i1.hdslb.com/bfs/openplatform/202407/VMRbFJZi1721029531102.jpeg"
$ws1.Range("B34").Value2 = "'2024-08-11"
$ws1.Range("C34").Value2 = "太仓·第六届龙狮动漫嘉年华"
$ws1.Range("D34").Value2 = "南园西路9号金仕堡4楼(南洋壹号公馆) AJ青少年篮球俱乐部"
$ws1.Range("E34").Value2 = "2024.08.11 09:00-08.11 17:00"
$ws1.Range("F34").Value2 = 29
$ws1.Range("G34").Value2 = 50
$ws1.Range("H34").Value2 = "https://show.bilibili.com/platform/detail.html?id=88517"
$ws1.Range("I34").Value2 = "//i0.hdslb.com/bfs/openplatform/202406/xNN6ZUtC1719579739903.jpeg"
$ws1.Range("C35").Value2 = "昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典"
$ws1.Range("D35").Value2 = "玉山镇前进西路1066号  华东海鲜城"
$ws1.Range("E35").Value2 = "2024.08.11 09:00-08.11 16:00"
$ws1.Range("F35").Value2 = 98
$ws1.Range("G35").Value2 = 45
$ws1.Range("H35").Value2 = "https://show.bilibili.com/platform/detail.html?id=87750"
$ws1.Range("I35").Value2 = "//i1.hdslb.com/bfs/openplatform/202406/rDcLl3X11719561046839.jpeg"
$ws1.Range("C36").Value2 = "苏州·OCG国潮动漫游戏嘉年华凌飞内场"
$ws1.Range("D36").Value2 = "苏州大道东688号 苏州国际博览中心"
$ws1.Range("E36").Value2 = "2024.08.11 09:00-08.11 17:00"
$ws1.Range("F36").Value2 = 27
$ws1.Range("G36").Value2 = 238
$ws1.Range("H36").Value2 = "https://show.bilibili.com/platform/detail.html?id=89239"
$ws1.Range("I36").Value2 = "//i0.hdslb.com/bfs/openplatform/202407/x77hdkJC1720332496348.jpeg"

# Sheet2 (演出)
$ws2.Range("F2").Value2 = 318

# Sheet4 (全部类型) simple F-value updates
$ws4.Range("F2").Value2 = 20853
$ws4.Range("F3").Value2 = 810
$ws4.Range("F4").Value2 = 335
$ws4.Range("F5").Value2 = 1118
$ws4.Range("F6").Value2 = 22
$ws4.Range("F7").Value2 = 7796
$ws4.Range("F8").Value2 = 544
$ws4.Range("F9").Value2 = 20
$ws4.Range("F10").Value2 = 753
$ws4.Range("F11").Value2 = 296
$ws4.Range("F12").Value2 = 52
$ws4.Range("F13").Value2 = 180
$ws4.Range("F14").Value2 = 149
$ws4.Range("F15").Value2 = 25
$ws4.Range("F18").Value2 = 1353
$ws4.Range("F19").Value2 = 500
$ws4.Range("F20").Value2 = 78
$ws4.Range("F23").Value2 = 77
$ws4.Range("F24").Value2 = 81
$ws4.Range("F25").Value2 = 341
$ws4.Range("F26").Value2 = 1161
$ws4.Range("F28").Value2 = 33
$ws4.Range("F29").Value2 = 207
$ws4.Range("F30").Value2 = 318
$ws4.Range("F31").Value2 = 590
$ws4.Range("F33").Value2 = 114
$ws4.Range("F39").Value2 = 58
$ws4.Range("F40").Value2 = 12933
$ws4.Range("F41").Value2 = 1352
$ws4.Range("F42").Value2 = 110
$ws4.Range("F43").Value2 = 43
$ws4.Range("F44").Value2 = 65
$ws4.Range("F45").Value2 = 293
$ws4.Range("F46").Value2 = 407
$ws4.Range("F47").Value2 = 4036

# Sheet4 rows 34-38 rotation updates
$ws4.Range("C34").Value2 = "苏州·爱乐之城·经典电影作品音乐会"
$ws4.Range("D34").Value2 = "念珠街121号道前街与吉庆街路口距养育巷地铁站 苏州市会议中心"
$ws4.Range("E34").Value2 = "2024.08.10 19:30-08.10 21:00"
$ws4.Range("F34").Value2 = 36
$ws4.Range("G34").Value2 = 90
$ws4.Range("H34").Value2 = "https://show.bilibili.com/platform/detail.html?id=86194"
$ws4.Range("I34").Value2 = "//i2.hdslb.com/bfs/openplatform/202405/vagzbfox1716438290025.jpeg"
$ws4.Range("C35").Value2 = "苏州·第三届.OCG.Summer Carnival-国潮动漫游戏嘉年华"
$ws4.Range("D35").Value2 = "苏州大道东688号 苏州国际博览中心"
$ws4.Range("E35").Value2 = "2024.08.10 09:00-08.11 17:00"
$ws4.Range("F35").Value2 = 4970
$ws4.Range("G35").Value2 = 70
$ws4.Range("H35").Value2 = "https://show.bilibili.com/platform/detail.html?id=88451"
$ws4.Range("I35").Value2 = "//i1.hdslb.com/bfs/openplatform/202407/VMRbFJZi1721029531102.jpeg"
$ws4.Range("B36").Value2 = "'2024-08-11"
$ws4.Range("C36").Value2 = "太仓·第六届龙狮动漫嘉年华"
$ws4.Range("D36").Value2 = "南园西路9号金仕堡4楼(南洋壹号公馆) AJ青少年篮球俱乐部"
$ws4.Range("E36").Value2 = "2024.08.11 09:00-08.11 17:00"
$ws4.Range("F36").Value2 = 29
$ws4.Range("G36").Value2 = 50
$ws4.Range("H36").Value2 = "https://show.bilibili.com/platform/detail.html?id=88517"
$ws4.Range("I36").Value2 = "//i0.hdslb.com/bfs/openplatform/202406/xNN6ZUtC1719579739903.jpeg"
$ws4.Range("C37").Value2 = "昆山·第七届·xcy新次元动漫嘉年华-狂欢盛典"
$ws4.Range("D37").Value2 = "玉山镇前进西路1066号  华东海鲜城"
$ws4.Range("E37").Value2 = "2024.08.11 09:00-08.11 16:00"
$ws4.Range("F37").Value2 = 98
$ws4.Range("G37").Value2 = 45
$ws4.Range("H37").Value2 = "https://show.bilibili.com/platform/detail.html?id=87750"
$ws4.Range("I37").Value2 = "//i1.hdslb.com/bfs/openplatform/202406/rDcLl3X11719561046839.jpeg"
$ws4.Range("C38").Value2 = "苏州·OCG国潮动漫游戏嘉年华凌飞内场"
$ws4.Range("D38").Value2 = "苏州大道东688号 苏州国际博览中心"
$ws4.Range("E38").Value2 = "2024.08.11 09:00-08.11 17:00"
$ws4.Range("F38").Value2 = 27
$ws4.Range("G38").Value2 = 238
$ws4.Range("H38").Value2 = "https://show.bilibili.com/platform/detail.html?id=89239"
$ws4.Range("I38").Value2 = "//i0.hdslb.com/bfs/openplatform/202407/x77hdkJC1720332496348.jpeg"
